$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row of data describing the "interact" dialogue entry
$ws.Range("A7").Value = 10005
$ws.Range("B7").Value = "interact"
$ws.Range("C7").Value = "玩家交互"
$ws.Range("D7").Value = "icon/down"
$ws.Range("E7").Value = "icon/up"
$ws.Range("F7").Value = "icon/down"
$ws.Range("G7").Value = "none"
$ws.Range("H7").Value = "null"
$ws.Range("I7").Value = 10001

# Update the active selection to match the saved view state
$ws.Range("H13").Select()
